# Update the "In Class Demonstration" Diff-in-Diff regression table:
#   - Shorten the header/row labels (drop the " Diff-in-Diff" suffix)
#   - Replace the reported coefficients with the new regression run's values
#   - Replace the r2_adj row with the new numbers
#
# Cells whose new text is a plain numeric-looking string (e.g. "-0.042",
# with no significance stars) must be written with a leading apostrophe via
# .Formula so Excel stores them as literal text (matching the source file,
# where these are shared strings) instead of silently re-parsing them as
# numbers. Cells that already contain non-numeric characters (the "*"
# significance markers) are safe to set directly with .Value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - column headers
$ws.Range("B1").Value = "C"
$ws.Range("C1").Value = "LF"
$ws.Range("D1").Value = "FFR"
$ws.Range("E1").Value = "A"

# Row 2 - "C" row
$ws.Range("A2").Value = "C"
$ws.Range("C2").Formula = "'-0.042"
$ws.Range("D2").Formula = "'-0.044"
$ws.Range("E2").Formula = "'-0.001"

# Row 3 - "LF" row
$ws.Range("A3").Value = "LF"
$ws.Range("B3").Formula = "'-0.368"
$ws.Range("D3").Value = "0.998***"
$ws.Range("E3").Value = "-0.016*"

# Row 4 - "FFR" row
$ws.Range("A4").Value = "FFR"
$ws.Range("B4").Formula = "'-0.158"
$ws.Range("C4").Value = "0.408***"
$ws.Range("E4").Formula = "'-0.004"

# Row 5 - "A" row
$ws.Range("A5").Value = "A"
$ws.Range("B5").Formula = "'-3.273"
$ws.Range("C5").Value = "-9.409*"
$ws.Range("D5").Formula = "'-5.344"

# Row 6 - "Constant" row (label unchanged)
$ws.Range("B6").Formula = "'-0.685"
$ws.Range("C6").Value = "-0.431*"
$ws.Range("D6").Formula = "'-0.429"
$ws.Range("E6").Value = "-0.021**"

# Row 7 - "r2_adj" row (label unchanged); these are genuine numbers
$ws.Range("B7").Value = -0.01
$ws.Range("C7").Value = 0.63
$ws.Range("D7").Value = 0.57
$ws.Range("E7").Value = 0.35
